# fix login and register page
# Append new username/password rows (rows 4-13) to the "Database" sheet,
# extending the used range from A1:B3 to A1:B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (username, password). Row 12 has no username (column A
# left blank), matching the source data.
$rows = @(
    @{ Row = 4;  User = "thanh";   Pass = "1234" },
    @{ Row = 5;  User = "hello";   Pass = "12345" },
    @{ Row = 6;  User = "thanh2";  Pass = "123" },
    @{ Row = 7;  User = "thanh3";  Pass = "123" },
    @{ Row = 8;  User = "thanh6";  Pass = "123" },
    @{ Row = 9;  User = "thanh7";  Pass = "1234" },
    @{ Row = 10; User = "hello1";  Pass = "123" },
    @{ Row = 11; User = "hello23"; Pass = "123" },
    @{ Row = 12; User = $null;     Pass = "1234" },
    @{ Row = 13; User = "thanh1";  Pass = "1234" }
)

# Passwords are numeric-looking strings that must stay text (matching the
# existing B2 = "1234" cell), so format column B as text before writing.
$ws.Range("B4:B13").NumberFormat = "@"

foreach ($r in $rows) {
    $rowNum = $r.Row
    if ($r.User -ne $null) {
        $ws.Cells.Item($rowNum, 1).Value = $r.User
    }
    $ws.Cells.Item($rowNum, 2).Value = $r.Pass
}
